# LOT2007.xlsx edit: the underlying source data shifted by one row.
# Row 13 (which held only the "427823 - Adriane Maria Ferreira Milagres"
# teacher-name values in B/C with no label in A) is removed, shifting
# everything below it up by one row. That also cascades into a second,
# independent one-row shift further down (Método:/Critério:/Avaliação
# block), which we fix up by writing the correct text into B/C directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the stray row (old row 13) - this shifts rows 14-25 up to 13-24
#    and Excel keeps row heights / cell formatting attached to the rows
#    that move, so the ht="60"/"120"/"30" customHeight attributes tag
#    along automatically.
$ws.Rows.Item(13).Delete()

# 2) Patch the cell values that don't line up after the straight shift.

# Row 10 (Objetivos:) B/C now needs the teacher name instead of the
# long objectives paragraph.
$ws.Range("B10").Value = "427823 - Adriane Maria Ferreira Milagres"
$ws.Range("C10").Value = "427823 - Adriane Maria Ferreira Milagres"

# Row 13 (Programa resumido:) B/C becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) B/C becomes "01/01/2018".
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18 (Método:) B/C becomes the teacher name again.
$ws.Range("B18").Value = "427823 - Adriane Maria Ferreira Milagres"
$ws.Range("C18").Value = "427823 - Adriane Maria Ferreira Milagres"

# Row 19 (Critério:) B/C becomes the evaluation-method sentence.
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

# Row 20 (Norma de recuperação:) B/C becomes the NF formula sentence.
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3"

# Row 21 (Bibliografia:) B/C becomes the recuperação (MR) sentence.
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
